# The "Year of Treatment" column (B) is removed entirely; all columns to
# its right (C:K -> B:J) shift one position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(2).Delete()
